$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook gets a new weekly record. The existing row 8 (old record,
# 2022-08-09) is pushed down to row 9 unchanged, and row 8 is overwritten
# with the new weekly record (2022-08-25).

# Duplicate row 8 (with all its values/formatting) down into a freshly
# inserted row 9.
$ws.Rows.Item(8).Copy()
$ws.Rows.Item(9).Insert()

# Now overwrite row 8 in place with the new record's values.
$ws.Range("D8").Value = 44798
$ws.Range("J8").Value = 80
$ws.Range("K8").Value = 7000
$ws.Range("L8").Value = 7000
$ws.Range("M8").Value = 7000
$ws.Range("O8").Value = "Provincia de Diguillín"
$ws.Range("P8").Value = 700
